# Updates the cryptos list table (Sheet1) with the latest scrape values
# (commit: "Updated cryptos list ... with GitHub Actions").
#
# The source workbook stores every data cell as literal text (prices keep
# "1.000s" separators, % cells keep their padding spaces), so every write
# below goes through as Text: we pre-format the target cell as Text ("@")
# before assigning, which stops Excel's COM layer from "helpfully"
# re-interpreting a value like "0.997" or "6.30" as a Number (and losing
# the trailing zero) instead of leaving it as the literal string we want.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.122.79"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.812.33"
$ws.Range("E3").Value = "  +1.39%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.82"
$ws.Range("E5").Value = "  -0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.07"
$ws.Range("E6").Value = "  -3.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.811.10"
$ws.Range("E7").Value = "  +1.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.997"
$ws.Range("E8").Value = "  -0.34%  "
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("E10").Value = "  +1.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.30"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.28"
$ws.Range("E13").Value = "  -3.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000247"
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.446.81"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.811.47"
$ws.Range("E16").Value = "  +1.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.257.01"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.44"
$ws.Range("E18").Value = "  +1.69%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.114"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.34"
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.48"
$ws.Range("E21").Value = "  +5.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "488.53"
$ws.Range("E22").Value = "  -1.59%  "
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000159"
$ws.Range("E24").Value = "  +3.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.85"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.27"
$ws.Range("E26").Value = "  -3.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.25"
$ws.Range("E27").Value = "  -1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.04"
$ws.Range("E28").Value = "  -2.46%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("E31").Value = "  -0.22%  "
$ws.Range("E32").Value = "  -4.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.959.97"
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.93"
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.757.72"
$ws.Range("E35").Value = "  +1.62%  "
$ws.Range("E36").Value = "  -1.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.141"
$ws.Range("E37").Value = "  +5.28%  "
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.92"
$ws.Range("E39").Value = "  +0.41%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.04"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "48.59"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "426.62"
$ws.Range("E44").Value = "  -2.83%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.99"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.834.91"
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.29"
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "39.55"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("E51").Value = "  -1.22%  "
